$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections on existing rows ---

# Row 3: CageID (col D) corrected from "12A" to "גד"
$ws.Range("D3").Value = "גד"

# Row 5: CageID (col D) corrected from string "50A" to a plain number 32
$ws.Range("D5").Value = 32

# --- New bird record appended as row 18 ---
$ws.Range("A18").Value = 8
$ws.Range("B18").Value = "American Gouldian"
$ws.Range("C18").Value = "North America"
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = "Female"
$ws.Range("F18").Value = 318692993
$ws.Range("G18").Value = 435345
$ws.Range("H18").Value = "15/05/2023"
$ws.Range("I18").Value = "Red"
$ws.Range("J18").Value = "Purple"
$ws.Range("K18").Value = "Green Pastel"
